$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2211221122112211
$ws.Range("C2").Value = 0.5016501650165016
$ws.Range("J2").Value = 0.0132013201320132
$ws.Range("O2").Value = 0.0033003300330033
$ws.Range("P2").Value = 0.1617161716171617
$ws.Range("S2").Value = 0.09900990099009901
$ws.Range("B3").Value = 0.006211180124223602
$ws.Range("C3").Value = 0.03726708074534162
$ws.Range("J3").Value = 0.06211180124223602
$ws.Range("P3").Value = 0.6583850931677019
$ws.Range("S3").Value = 0.2360248447204969
$ws.Range("J4").Value = 0.06451612903225806
$ws.Range("P4").Value = 0.6451612903225806
$ws.Range("S4").Value = 0.2903225806451613
$ws.Range("B6").Value = 0.04296875
$ws.Range("D6").Value = 0.00390625
$ws.Range("E6").Value = 0.00390625
$ws.Range("F6").Value = 0.078125
$ws.Range("J6").Value = 0.2890625
$ws.Range("O6").Value = 0.03515625
$ws.Range("Q6").Value = 0.12109375
$ws.Range("R6").Value = 0.0546875
$ws.Range("S6").Value = 0.37109375
$ws.Range("B7").Value = 0.1382113821138211
$ws.Range("D7").Value = 0.01219512195121951
$ws.Range("F7").Value = 0.06097560975609756
$ws.Range("J7").Value = 0.1016260162601626
$ws.Range("O7").Value = 0.01626016260162602
$ws.Range("Q7").Value = 0.1991869918699187
$ws.Range("R7").Value = 0.08130081300813008
$ws.Range("S7").Value = 0.3902439024390244
$ws.Range("B8").Value = 0.09611451942740286
$ws.Range("D8").Value = 0.016359918200409
$ws.Range("F8").Value = 0.08588957055214724
$ws.Range("J8").Value = 0.09406952965235174
$ws.Range("O8").Value = 0.03476482617586912
$ws.Range("Q8").Value = 0.1513292433537832
$ws.Range("R8").Value = 0.1165644171779141
$ws.Range("S8").Value = 0.4049079754601227
$ws.Range("B9").Value = 0.1468926553672316
$ws.Range("D9").Value = 0.005649717514124294
$ws.Range("F9").Value = 0.03389830508474576
$ws.Range("J9").Value = 0.1016949152542373
$ws.Range("O9").Value = 0.02824858757062147
$ws.Range("Q9").Value = 0.1638418079096045
$ws.Range("R9").Value = 0.07344632768361582
$ws.Range("S9").Value = 0.4463276836158192
$ws.Range("B10").Value = 0.09809027777777778
$ws.Range("D10").Value = 0.01475694444444444
$ws.Range("E10").Value = 0.0008680555555555555
$ws.Range("F10").Value = 0.07465277777777778
$ws.Range("J10").Value = 0.1189236111111111
$ws.Range("O10").Value = 0.03298611111111111
$ws.Range("Q10").Value = 0.1866319444444444
$ws.Range("R10").Value = 0.08940972222222222
$ws.Range("S10").Value = 0.3836805555555556
$ws.Range("G11").Value = 0.1454081632653061
$ws.Range("J11").Value = 0.09693877551020408
$ws.Range("K11").Value = 0.2168367346938775
$ws.Range("L11").Value = 0.5306122448979592
$ws.Range("S11").Value = 0.01020408163265306
$ws.Range("G12").Value = 0.7442922374429224
$ws.Range("J12").Value = 0.1872146118721461
$ws.Range("K12").Value = 0.0136986301369863
$ws.Range("L12").Value = 0.0410958904109589
$ws.Range("S12").Value = 0.0136986301369863
$ws.Range("F15").Value = 0.02962962962962963
$ws.Range("H15").Value = 0.1814814814814815
$ws.Range("I15").Value = 0.04074074074074074
$ws.Range("J15").Value = 0.2777777777777778
$ws.Range("K15").Value = 0.07037037037037037
$ws.Range("M15").Value = 0.01851851851851852
$ws.Range("N15").Value = 0.003703703703703704
$ws.Range("O15").Value = 0.05925925925925926
$ws.Range("S15").Value = 0.3185185185185185
$ws.Range("F16").Value = 0.01176470588235294
$ws.Range("I16").Value = 0.07647058823529412
$ws.Range("J16").Value = 0.3470588235294118
$ws.Range("K16").Value = 0.1411764705882353
$ws.Range("M16").Value = 0.02941176470588235
$ws.Range("O16").Value = 0.06470588235294118
$ws.Range("S16").Value = 0.1294117647058824
$ws.Range("F17").Value = 0.02557544757033248
$ws.Range("H17").Value = 0.2046035805626598
$ws.Range("I17").Value = 0.07672634271099744
$ws.Range("J17").Value = 0.3350383631713555
$ws.Range("K17").Value = 0.1585677749360614
$ws.Range("M17").Value = 0.01534526854219949
$ws.Range("N17").Value = 0.002557544757033248
$ws.Range("O17").Value = 0.07928388746803069
$ws.Range("S17").Value = 0.1023017902813299
$ws.Range("F18").Value = 0.01463414634146342
$ws.Range("H18").Value = 0.2097560975609756
$ws.Range("I18").Value = 0.1024390243902439
$ws.Range("J18").Value = 0.3365853658536586
$ws.Range("K18").Value = 0.0975609756097561
$ws.Range("M18").Value = 0.02926829268292683
$ws.Range("O18").Value = 0.06341463414634146
$ws.Range("S18").Value = 0.1463414634146341
$ws.Range("F19").Value = 0.02925327174749807
$ws.Range("H19").Value = 0.2193995381062356
$ws.Range("I19").Value = 0.07852193995381063
$ws.Range("J19").Value = 0.3302540415704388
$ws.Range("K19").Value = 0.1339491916859122
$ws.Range("M19").Value = 0.02463433410315627
$ws.Range("N19").Value = 0.0007698229407236335
$ws.Range("O19").Value = 0.07775211701308699
$ws.Range("S19").Value = 0.1054657428791378
